$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: The "Make sure the schema name is not conflicting ..." notice
# (together with the blank paragraph that followed it) is removed from
# just after the application.properties intro line -- that sentence
# reappears, reworded, as a new italic "NB:" paragraph later on (see
# Part 2 below). Removing it here simply lets the flyway config block
# (url/user/password/enabled) slide up into its old place.
# -----------------------------------------------------------------------

$findRange = $d.Content
$findRange.Find.Execute("Make sure the schema name is not conflicting with another one of your schemas, if so, please change it.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$schemaParaIndex = $findRange.Paragraphs.Item(1).Index

$startDel = $d.Paragraphs.Item($schemaParaIndex).Range.Start
$endDel = $d.Paragraphs.Item($schemaParaIndex + 1).Range.End
$delRange = $d.Range($startDel, $endDel)
$delRange.Delete()

# -----------------------------------------------------------------------
# Part 2: Add a new italic paragraph right after the datasource password
# paragraph ("spring.datasource.password= your_password_here"), made up
# of two runs:
#   "NB: Make sure the schema name is not conflicting with another one
#    of your schemas, if so, please change it."
#   " As it is in the V1 of flyway migration scripts it will also have
#    to be changed there."
# -----------------------------------------------------------------------

$findRange2 = $d.Content
$findRange2.Find.Execute("spring.datasource.password=", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$passwordParaIndex = $findRange2.Paragraphs.Item(1).Index

$passwordPara = $d.Paragraphs.Item($passwordParaIndex)
$insertionPoint = $passwordPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($passwordParaIndex + 1)
$newPara.Range.Text = "NB: Make sure the schema name is not conflicting with another one of your schemas, if so, please change it."

$startPos = $newPara.Range.Start
$midPos = $newPara.Range.End - 1

$secondRunRange = $newPara.Range.Duplicate
$secondRunRange.SetRange($midPos, $midPos)
$secondRunRange.InsertAfter(" As it is in the V1 of flyway migration scripts it will also have to be changed there.")

$endPos = $newPara.Range.End - 1

$firstRunRange = $newPara.Range.Duplicate
$firstRunRange.SetRange($startPos, $midPos)
$firstRunRange.Font.Italic = 1
$firstRunRange.Font.ItalicBi = 1

$secondRunRange2 = $newPara.Range.Duplicate
$secondRunRange2.SetRange($midPos, $endPos)
$secondRunRange2.Font.Italic = 1
$secondRunRange2.Font.ItalicBi = 1

$fullNewParaRange = $newPara.Range
$fullNewParaRange.Font.Italic = 1
$fullNewParaRange.Font.ItalicBi = 1
